# Clear the per-site/building breakdown columns (AB:AK) and the
# "DIFFERENCE" column (AM) for all data rows, keeping the
# "PREVIOUS ACCOMPLISHMENT" column (AL) intact, matching the most
# updated status/accomplishment data as of May.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

$ws.Range("AB2:AK$lastRow").ClearContents()
$ws.Range("AM2:AM$lastRow").ClearContents()
